$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -21.04502370946747
$ws.Cells.Item(2, 3).Value = 1.913546882369055
$ws.Cells.Item(2, 4).Value = -21.04502370946747
$ws.Cells.Item(2, 5).Value = -21.04502370946747
$ws.Cells.Item(2, 6).Value = -21.04502370946747
$ws.Cells.Item(2, 7).Value = -21.04502370946747
$ws.Cells.Item(2, 8).Value = -21.04502370946747
$ws.Cells.Item(2, 9).Value = -21.04502370946747
$ws.Cells.Item(2, 10).Value = -21.04502370946747
$ws.Cells.Item(2, 11).Value = -21.04502370946747
$ws.Cells.Item(3, 2).Value = -21.04502370946747
$ws.Cells.Item(3, 3).Value = -21.04502370946747
$ws.Cells.Item(3, 4).Value = -21.04502370946747
$ws.Cells.Item(3, 5).Value = -21.04502370946747
$ws.Cells.Item(3, 6).Value = -21.04502370946747
$ws.Cells.Item(3, 7).Value = -21.04502370946747
$ws.Cells.Item(3, 8).Value = -21.04502370946747
$ws.Cells.Item(3, 9).Value = 1.244682550136331
$ws.Cells.Item(3, 10).Value = -21.04502370946747
$ws.Cells.Item(3, 11).Value = -21.04502370946747
$ws.Cells.Item(4, 2).Value = -21.04502370946747
$ws.Cells.Item(4, 3).Value = 2.014779027317557
$ws.Cells.Item(4, 4).Value = 1.635087886198406
$ws.Cells.Item(4, 5).Value = -21.04502370946747
$ws.Cells.Item(4, 6).Value = 3.412323162142619
$ws.Cells.Item(4, 7).Value = -21.04502370946747
$ws.Cells.Item(4, 8).Value = 1.492138988912234
$ws.Cells.Item(4, 9).Value = -21.04502370946747
$ws.Cells.Item(4, 10).Value = 0.900529075375736
$ws.Cells.Item(4, 11).Value = -21.04502370946747
$ws.Cells.Item(5, 2).Value = -21.04502370946747
$ws.Cells.Item(5, 3).Value = 1.674353621437287
$ws.Cells.Item(5, 4).Value = -21.04502370946747
$ws.Cells.Item(5, 5).Value = -21.04502370946747
$ws.Cells.Item(5, 6).Value = -21.04502370946747
$ws.Cells.Item(5, 7).Value = 2.831043617566067
$ws.Cells.Item(5, 8).Value = -21.04502370946747
$ws.Cells.Item(5, 9).Value = -21.04502370946747
$ws.Cells.Item(5, 10).Value = -21.04502370946747
$ws.Cells.Item(5, 11).Value = -21.04502370946747
$ws.Cells.Item(6, 2).Value = -21.04502370946747
$ws.Cells.Item(6, 3).Value = -21.04502370946747
$ws.Cells.Item(6, 4).Value = -21.04502370946747
$ws.Cells.Item(6, 5).Value = -21.04502370946747
$ws.Cells.Item(6, 6).Value = -21.04502370946747
$ws.Cells.Item(6, 7).Value = -21.04502370946747
$ws.Cells.Item(6, 8).Value = -21.04502370946747
$ws.Cells.Item(6, 9).Value = -21.04502370946747
$ws.Cells.Item(6, 10).Value = -21.04502370946747
$ws.Cells.Item(6, 11).Value = -21.04502370946747
$ws.Cells.Item(7, 2).Value = -21.04502370946747
$ws.Cells.Item(7, 3).Value = -21.04502370946747
$ws.Cells.Item(7, 4).Value = -21.04502370946747
$ws.Cells.Item(7, 5).Value = -21.04502370946747
$ws.Cells.Item(7, 6).Value = -21.04502370946747
$ws.Cells.Item(7, 7).Value = -21.04502370946747
$ws.Cells.Item(7, 8).Value = -21.04502370946747
$ws.Cells.Item(7, 9).Value = -21.04502370946747
$ws.Cells.Item(7, 10).Value = -21.04502370946747
$ws.Cells.Item(7, 11).Value = -21.04502370946747
$ws.Cells.Item(8, 2).Value = -21.04502370946747
$ws.Cells.Item(8, 3).Value = -21.04502370946747
$ws.Cells.Item(8, 4).Value = -21.04502370946747
$ws.Cells.Item(8, 5).Value = 1.89003574208616
$ws.Cells.Item(8, 6).Value = -21.04502370946747
$ws.Cells.Item(8, 7).Value = -21.04502370946747
$ws.Cells.Item(8, 8).Value = -21.04502370946747
$ws.Cells.Item(8, 9).Value = -21.04502370946747
$ws.Cells.Item(8, 10).Value = -21.04502370946747
$ws.Cells.Item(8, 11).Value = -21.04502370946747
$ws.Cells.Item(9, 2).Value = 4.321927461433711
$ws.Cells.Item(9, 3).Value = -21.04502370946747
$ws.Cells.Item(9, 4).Value = -21.04502370946747
$ws.Cells.Item(9, 5).Value = -21.04502370946747
$ws.Cells.Item(9, 6).Value = -21.04502370946747
$ws.Cells.Item(9, 7).Value = -21.04502370946747
$ws.Cells.Item(9, 8).Value = -21.04502370946747
$ws.Cells.Item(9, 9).Value = -21.04502370946747
$ws.Cells.Item(9, 10).Value = -21.04502370946747
$ws.Cells.Item(9, 11).Value = -21.04502370946747
$ws.Cells.Item(10, 2).Value = -21.04502370946747
$ws.Cells.Item(10, 3).Value = -21.04502370946747
$ws.Cells.Item(10, 4).Value = -21.04502370946747
$ws.Cells.Item(10, 5).Value = -21.04502370946747
$ws.Cells.Item(10, 6).Value = -21.04502370946747
$ws.Cells.Item(10, 7).Value = -21.04502370946747
$ws.Cells.Item(10, 8).Value = -21.04502370946747
$ws.Cells.Item(10, 9).Value = 1.739259072338638
$ws.Cells.Item(10, 10).Value = -21.04502370946747
$ws.Cells.Item(10, 11).Value = 2.222211185095067
$ws.Cells.Item(11, 2).Value = -21.04502370946747
$ws.Cells.Item(11, 3).Value = -21.04502370946747
$ws.Cells.Item(11, 4).Value = -21.04502370946747
$ws.Cells.Item(11, 5).Value = 2.889747781723523
$ws.Cells.Item(11, 6).Value = -21.04502370946747
$ws.Cells.Item(11, 7).Value = 2.835224287239107
$ws.Cells.Item(11, 8).Value = -21.04502370946747
$ws.Cells.Item(11, 9).Value = -21.04502370946747
$ws.Cells.Item(11, 10).Value = -21.04502370946747
$ws.Cells.Item(11, 11).Value = 1.95073464449619
$ws.Cells.Item(12, 2).Value = -21.04502370946747
$ws.Cells.Item(12, 3).Value = -21.04502370946747
$ws.Cells.Item(12, 4).Value = -21.04502370946747
$ws.Cells.Item(12, 5).Value = -21.04502370946747
$ws.Cells.Item(12, 6).Value = -21.04502370946747
$ws.Cells.Item(12, 7).Value = -21.04502370946747
$ws.Cells.Item(12, 8).Value = -21.04502370946747
$ws.Cells.Item(12, 9).Value = -21.04502370946747
$ws.Cells.Item(12, 10).Value = -21.04502370946747
$ws.Cells.Item(12, 11).Value = -21.04502370946747
$ws.Cells.Item(13, 2).Value = -21.04502370946747
$ws.Cells.Item(13, 3).Value = -21.04502370946747
$ws.Cells.Item(13, 4).Value = -21.04502370946747
$ws.Cells.Item(13, 5).Value = 2.518799230071247
$ws.Cells.Item(13, 6).Value = -21.04502370946747
$ws.Cells.Item(13, 7).Value = -21.04502370946747
$ws.Cells.Item(13, 8).Value = -21.04502370946747
$ws.Cells.Item(13, 9).Value = -21.04502370946747
$ws.Cells.Item(13, 10).Value = 1.683898439579075
$ws.Cells.Item(13, 11).Value = 1.770822558828672
$ws.Cells.Item(14, 2).Value = -21.04502370946747
$ws.Cells.Item(14, 3).Value = -21.04502370946747
$ws.Cells.Item(14, 4).Value = 1.535903024942266
$ws.Cells.Item(14, 5).Value = -21.04502370946747
$ws.Cells.Item(14, 6).Value = -21.04502370946747
$ws.Cells.Item(14, 7).Value = -21.04502370946747
$ws.Cells.Item(14, 8).Value = -21.04502370946747
$ws.Cells.Item(14, 9).Value = -21.04502370946747
$ws.Cells.Item(14, 10).Value = -21.04502370946747
$ws.Cells.Item(14, 11).Value = 1.947329144454159
$ws.Cells.Item(15, 2).Value = -21.04502370946747
$ws.Cells.Item(15, 3).Value = -21.04502370946747
$ws.Cells.Item(15, 4).Value = 1.74712255667472
$ws.Cells.Item(15, 5).Value = -21.04502370946747
$ws.Cells.Item(15, 6).Value = -21.04502370946747
$ws.Cells.Item(15, 7).Value = -21.04502370946747
$ws.Cells.Item(15, 8).Value = -21.04502370946747
$ws.Cells.Item(15, 9).Value = -21.04502370946747
$ws.Cells.Item(15, 10).Value = -21.04502370946747
$ws.Cells.Item(15, 11).Value = -21.04502370946747
$ws.Cells.Item(16, 2).Value = -21.04502370946747
$ws.Cells.Item(16, 3).Value = -21.04502370946747
$ws.Cells.Item(16, 4).Value = -21.04502370946747
$ws.Cells.Item(16, 5).Value = -21.04502370946747
$ws.Cells.Item(16, 6).Value = -21.04502370946747
$ws.Cells.Item(16, 7).Value = -21.04502370946747
$ws.Cells.Item(16, 8).Value = -21.04502370946747
$ws.Cells.Item(16, 9).Value = -21.04502370946747
$ws.Cells.Item(16, 10).Value = 1.926725584100465
$ws.Cells.Item(16, 11).Value = -21.04502370946747
$ws.Cells.Item(17, 2).Value = -21.04502370946747
$ws.Cells.Item(17, 3).Value = 2.126996734827298
$ws.Cells.Item(17, 4).Value = 1.850928889435287
$ws.Cells.Item(17, 5).Value = -21.04502370946747
$ws.Cells.Item(17, 6).Value = -21.04502370946747
$ws.Cells.Item(17, 7).Value = -21.04502370946747
$ws.Cells.Item(17, 8).Value = 2.059720435683475
$ws.Cells.Item(17, 9).Value = 2.096930762846415
$ws.Cells.Item(17, 10).Value = 2.526031845847624
$ws.Cells.Item(17, 11).Value = -21.04502370946747
$ws.Cells.Item(18, 2).Value = -21.04502370946747
$ws.Cells.Item(18, 3).Value = -21.04502370946747
$ws.Cells.Item(18, 4).Value = -21.04502370946747
$ws.Cells.Item(18, 5).Value = -21.04502370946747
$ws.Cells.Item(18, 6).Value = -21.04502370946747
$ws.Cells.Item(18, 7).Value = -21.04502370946747
$ws.Cells.Item(18, 8).Value = 1.996457148781035
$ws.Cells.Item(18, 9).Value = 2.044940220623787
$ws.Cells.Item(18, 10).Value = 2.421836537672778
$ws.Cells.Item(18, 11).Value = -21.04502370946747
$ws.Cells.Item(19, 2).Value = -21.04502370946747
$ws.Cells.Item(19, 3).Value = -21.04502370946747
$ws.Cells.Item(19, 4).Value = 2.05192640036502
$ws.Cells.Item(19, 5).Value = -21.04502370946747
$ws.Cells.Item(19, 6).Value = -21.04502370946747
$ws.Cells.Item(19, 7).Value = -21.04502370946747
$ws.Cells.Item(19, 8).Value = 1.637777397078668
$ws.Cells.Item(19, 9).Value = 1.825184647364066
$ws.Cells.Item(19, 10).Value = -21.04502370946747
$ws.Cells.Item(19, 11).Value = -21.04502370946747
$ws.Cells.Item(20, 2).Value = -21.04502370946747
$ws.Cells.Item(20, 3).Value = 1.077513950647304
$ws.Cells.Item(20, 4).Value = 1.527785902421053
$ws.Cells.Item(20, 5).Value = -21.04502370946747
$ws.Cells.Item(20, 6).Value = 3.225486997290401
$ws.Cells.Item(20, 7).Value = -21.04502370946747
$ws.Cells.Item(20, 8).Value = 1.656113719072802
$ws.Cells.Item(20, 9).Value = 1.228741447107211
$ws.Cells.Item(20, 10).Value = -21.04502370946747
$ws.Cells.Item(20, 11).Value = 2.070085523988798
$ws.Cells.Item(21, 2).Value = -21.04502370946747
$ws.Cells.Item(21, 3).Value = 1.334129610307981
$ws.Cells.Item(21, 4).Value = -21.04502370946747
$ws.Cells.Item(21, 5).Value = 1.655851350137878
$ws.Cells.Item(21, 6).Value = -21.04502370946747
$ws.Cells.Item(21, 7).Value = 2.522951508237797
$ws.Cells.Item(21, 8).Value = 1.466844979390713
$ws.Cells.Item(21, 9).Value = -21.04502370946747
$ws.Cells.Item(21, 10).Value = -21.04502370946747
$ws.Cells.Item(21, 11).Value = -21.04502370946747
